$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text would otherwise be auto-parsed as a number by Excel;
# force them to keep a Text format so the literal string is preserved exactly.
$textCells = @("D5","D8","D9","D11","D15","D16","D18","D20","D22","D23","D24","D25","D26","D28","D31","D32","D35","D36","D37","D40","D42","D44","D48","D50")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.438.32"
$ws.Range("E2").Value = "  -0.54%  "

$ws.Range("D3").Value = "1.569.29"
$ws.Range("E3").Value = "  -1.28%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "207.43"
$ws.Range("E5").Value = "  +0.05%  "

$ws.Range("E6").Value = "  -0.97%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "21.99"
$ws.Range("E8").Value = "  -1.15%  "

$ws.Range("D9").Value = "0.248"
$ws.Range("E9").Value = "  -1.27%  "

$ws.Range("E10").Value = "  -0.25%  "

$ws.Range("D11").Value = "0.0867"
$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("D12").Value = "1.793.64"
$ws.Range("E12").Value = "  -1.32%  "

$ws.Range("D13").Value = "1.587.59"
$ws.Range("E13").Value = "  +0.01%  "

$ws.Range("E14").Value = "  -0.93%  "

$ws.Range("D15").Value = "0.518"
$ws.Range("E15").Value = "  -3.10%  "

$ws.Range("D16").Value = "63.18"
$ws.Range("E16").Value = "  +0.00%  "

$ws.Range("D17").Value = "27.448.33"
$ws.Range("E17").Value = "  -0.40%  "

$ws.Range("D18").Value = "213.54"
$ws.Range("E18").Value = "  -1.81%  "

$ws.Range("E19").Value = "  -0.46%  "

$ws.Range("D20").Value = "7.24"
$ws.Range("E20").Value = "  -1.49%  "

$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("D22").Value = "4.11"
$ws.Range("E22").Value = "  -1.13%  "

$ws.Range("D23").Value = "9.63"
$ws.Range("E23").Value = "  +0.56%  "

$ws.Range("D24").Value = "2.00"
$ws.Range("E24").Value = "  +0.75%  "

$ws.Range("D25").Value = "154.50"
$ws.Range("E25").Value = "  +0.64%  "

$ws.Range("D26").Value = "6.83"
$ws.Range("E26").Value = "  +1.12%  "

$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("D28").Value = "14.99"
$ws.Range("E28").Value = "  -0.63%  "

$ws.Range("E29").Value = "  -1.62%  "

$ws.Range("E30").Value = "  -0.76%  "

$ws.Range("D31").Value = "0.0470"
$ws.Range("E31").Value = "  +0.95%  "

$ws.Range("D32").Value = "3.19"
$ws.Range("E32").Value = "  -2.05%  "

$ws.Range("D33").Value = "1.359.53"
$ws.Range("E33").Value = "  -1.07%  "

$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("D35").Value = "1.53"
$ws.Range("E35").Value = "  +1.09%  "

$ws.Range("D36").Value = "0.967"
$ws.Range("E36").Value = "  +0.38%  "

$ws.Range("D37").Value = "2.30"
$ws.Range("E37").Value = "  -0.27%  "

$ws.Range("E38").Value = "  +1.08%  "

$ws.Range("E39").Value = "  -1.13%  "

$ws.Range("D40").Value = "0.818"
$ws.Range("E40").Value = "  +0.61%  "

$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("D42").Value = "0.973"
$ws.Range("E42").Value = "  +0.46%  "

$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("D44").Value = "64.06"
$ws.Range("E44").Value = "  +0.39%  "

$ws.Range("E45").Value = "  -0.91%  "

$ws.Range("E46").Value = "  -1.96%  "

$ws.Range("D47").Value = "1.706.64"
$ws.Range("E47").Value = "  -1.22%  "

$ws.Range("D48").Value = "85.09"
$ws.Range("E48").Value = "  -2.26%  "

$ws.Range("D49").Value = "0.0₇0995"
$ws.Range("E49").Value = "  +3.41%  "

$ws.Range("D50").Value = "0.0954"
$ws.Range("E50").Value = "  -1.46%  "

$ws.Range("E51").Value = "  -0.47%  "
